$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.93351310881776
$ws.Range("C2").Value = 11.64185014084704
$ws.Range("D2").Value = 3.918134423018945
$ws.Range("F2").Value = 15.88031944475604
$ws.Range("G2").Value = 14.15776795562885
$ws.Range("H2").Value = 10.7760138999018
$ws.Range("O2").Value = 14.42781711171567

$ws.Range("B3").Value = 14.06371179874951
$ws.Range("C3").Value = 11.08134734701999
$ws.Range("D3").Value = 3.778557731646563
$ws.Range("F3").Value = 16.00516197081918
$ws.Range("G3").Value = 14.3771886394656
$ws.Range("H3").Value = 10.85000455132577
$ws.Range("O3").Value = 14.56937377645897

$ws.Range("B4").Value = 13.49971208030673
$ws.Range("C4").Value = 10.72178393167878
$ws.Range("D4").Value = 3.689628288360711
$ws.Range("F4").Value = 16.089396163542
$ws.Range("G4").Value = 14.52335346772309
$ws.Range("H4").Value = 10.89799245352996
$ws.Range("O4").Value = 14.66161291441114

$ws.Range("B5").Value = 13.26245641121145
$ws.Range("C5").Value = 10.57153525370216
$ws.Range("D5").Value = 3.652614477717628
$ws.Range("F5").Value = 16.12561389044672
$ws.Range("G5").Value = 14.5857549813185
$ws.Range("H5").Value = 10.91819159840708
$ws.Range("O5").Value = 14.7005372909308

$ws.Range("B6").Value = 13.22261571549053
$ws.Range("C6").Value = 10.54636661731432
$ws.Range("D6").Value = 3.64642277562597
$ws.Range("F6").Value = 16.13174163935838
$ws.Range("G6").Value = 14.59628702441666
$ws.Range("H6").Value = 10.92158454684935
$ws.Range("O6").Value = 14.70708129262694

$ws.Range("B7").Value = 13.49654225501884
$ws.Range("C7").Value = 10.71977248305399
$ws.Range("D7").Value = 3.689132189699385
$ws.Range("F7").Value = 16.08987696865612
$ws.Range("G7").Value = 14.52418359460567
$ws.Range("H7").Value = 10.89826225916296
$ws.Range("O7").Value = 14.66213245491905

$ws.Range("B8").Value = 14.63987572820987
$ws.Range("C8").Value = 11.45187115628012
$ws.Range("D8").Value = 3.870699756529795
$ws.Range("F8").Value = 15.92178311363235
$ws.Range("G8").Value = 14.23102640364129
$ws.Range("H8").Value = 10.80099547008109
$ws.Range("O8").Value = 14.4755197002829

$ws.Range("B9").Value = 16.64088930889133
$ws.Range("C9").Value = 12.75973928267894
$ws.Range("D9").Value = 4.199692340406275
$ws.Range("F9").Value = 15.65294512904592
$ws.Range("G9").Value = 13.74869725770876
$ws.Range("H9").Value = 10.63051852500798
$ws.Range("O9").Value = 14.15191305123695

$ws.Range("B10").Value = 17.96069348421627
$ws.Range("C10").Value = 13.63680996244116
$ws.Range("D10").Value = 4.423244813064995
$ws.Range("F10").Value = 15.49336899027122
$ws.Range("G10").Value = 13.45328739305439
$ws.Range("H10").Value = 10.51758212158674
$ws.Range("O10").Value = 13.94012524316598

$ws.Range("B11").Value = 18.52803864556936
$ws.Range("C11").Value = 14.01669900617182
$ws.Range("D11").Value = 4.52071995191422
$ws.Range("F11").Value = 15.42919723791028
$ws.Range("G11").Value = 13.33227278682441
$ws.Range("H11").Value = 10.46887032655298
$ws.Range("O11").Value = 13.84945330643957

$ws.Range("B12").Value = 18.73810404747452
$ws.Range("C12").Value = 14.15774891536712
$ws.Range("D12").Value = 4.557004920424767
$ws.Range("F12").Value = 15.40612293075563
$ws.Range("G12").Value = 13.28841938616425
$ws.Range("H12").Value = 10.45080702052996
$ws.Range("O12").Value = 13.81593742688793

$ws.Range("B13").Value = 18.69307541936708
$ws.Range("C13").Value = 14.12749692055636
$ws.Range("D13").Value = 4.549218466740735
$ws.Range("F13").Value = 15.41103762140832
$ws.Range("G13").Value = 13.29777548073127
$ws.Range("H13").Value = 10.45468025483112
$ws.Range("O13").Value = 13.82311915098202

$ws.Range("B14").Value = 18.54541668066731
$ws.Range("C14").Value = 14.02835980393505
$ws.Range("D14").Value = 4.523717803646143
$ws.Range("F14").Value = 15.42727424019457
$ws.Range("G14").Value = 13.32862512576782
$ws.Range("H14").Value = 10.46737657683493
$ws.Range("O14").Value = 13.84667948429079

$ws.Range("B15").Value = 18.45434888655165
$ws.Range("C15").Value = 13.96726846741076
$ws.Range("D15").Value = 4.508015744729603
$ws.Range("F15").Value = 15.43737976840836
$ws.Range("G15").Value = 13.34777976170918
$ws.Range("H15").Value = 10.47520328595374
$ws.Range("O15").Value = 13.86121773214982

$ws.Range("B16").Value = 17.92294897558831
$ws.Range("C16").Value = 13.61159369058145
$ws.Range("D16").Value = 4.416787900831846
$ws.Range("F16").Value = 15.49773351740612
$ws.Range("G16").Value = 13.46146922443758
$ws.Range("H16").Value = 10.5208191292885
$ws.Range("O16").Value = 13.9461653329837

$ws.Range("B17").Value = 17.58846951518376
$ws.Range("C17").Value = 13.38846168261691
$ws.Range("D17").Value = 4.359726676577278
$ws.Range("F17").Value = 15.53692619118348
$ws.Range("G17").Value = 13.53467312744486
$ws.Range("H17").Value = 10.54948499671939
$ws.Range("O17").Value = 13.99973333430572

$ws.Range("B18").Value = 17.39297672391044
$ws.Range("C18").Value = 13.25832754926631
$ws.Range("D18").Value = 4.326510369335662
$ws.Range("F18").Value = 15.56026043026888
$ws.Range("G18").Value = 13.57803335296843
$ws.Range("H18").Value = 10.56622356082576
$ws.Range("O18").Value = 14.03107779135282

$ws.Range("B19").Value = 17.32625313749134
$ws.Range("C19").Value = 13.21396017859634
$ws.Range("D19").Value = 4.315196502799649
$ws.Range("F19").Value = 15.56829649838602
$ws.Range("G19").Value = 13.59292847885183
$ws.Range("H19").Value = 10.57193401912633
$ws.Range("O19").Value = 14.04178200535006

$ws.Range("B20").Value = 17.62439740603285
$ws.Range("C20").Value = 13.41240066108547
$ws.Range("D20").Value = 4.365842099489477
$ws.Range("F20").Value = 15.53267201958033
$ws.Range("G20").Value = 13.52675017388486
$ws.Range("H20").Value = 10.54640751983816
$ws.Range("O20").Value = 13.99397567900035

$ws.Range("B21").Value = 18.58891740364689
$ws.Range("C21").Value = 14.05755535021541
$ws.Range("D21").Value = 4.531225121760395
$ws.Range("F21").Value = 15.42247175367673
$ws.Range("G21").Value = 13.31950990937622
$ws.Range("H21").Value = 10.46363697191157
$ws.Range("O21").Value = 13.83973696541236

$ws.Range("B22").Value = 19.19144350231936
$ws.Range("C22").Value = 14.462826463798
$ws.Range("D22").Value = 4.635651910585806
$ws.Range("F22").Value = 15.35760332284283
$ws.Range("G22").Value = 13.19558585901097
$ws.Range("H22").Value = 10.41177277978981
$ws.Range("O22").Value = 13.74371247017933

$ws.Range("B23").Value = 18.8724177421328
$ws.Range("C23").Value = 14.24804093782735
$ws.Range("D23").Value = 4.580258206871737
$ws.Range("F23").Value = 15.39156525766298
$ws.Range("G23").Value = 13.26065567202994
$ws.Range("H23").Value = 10.43924958562871
$ws.Range("O23").Value = 13.7945238193221

$ws.Range("B24").Value = 17.60816435976486
$ws.Range("C24").Value = 13.40158361462677
$ws.Range("D24").Value = 4.363078595584036
$ws.Range("F24").Value = 15.53459283403717
$ws.Range("G24").Value = 13.53032817553294
$ws.Range("H24").Value = 10.5477980429121
$ws.Range("O24").Value = 13.99657700983326

$ws.Range("B25").Value = 16.12577964503738
$ws.Range("C25").Value = 12.42028705777611
$ws.Range("D25").Value = 4.113773622268636
$ws.Range("F25").Value = 15.71907350869982
$ws.Range("G25").Value = 13.8690155013069
$ws.Range("H25").Value = 10.67447204459421
$ws.Range("O25").Value = 14.23490930322805

